# Hortaliza, Feria Lagunitas de Puerto Montt - Pepino dulce
# Insert two new weekly records (rows 53-54) above the existing data,
# pushing the prior rows 53-103 down to 55-105.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 53:54 - this shifts old rows 53-103 down to 55-105
# and extends the used range to A1:R105 automatically.
$ws.Rows("53:54").Insert()

# New row 53: Especial quality entry dated 2023-06-02 (serial 45079)
$ws.Range("A53").Value = 4
$ws.Range("B53").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C53").Value = "Los Lagos"
$ws.Range("D53").Value = 45079
$ws.Range("E53").Value = 10
$ws.Range("F53").Value = 100112043
$ws.Range("G53").Value = "Pepino dulce"
$ws.Range("H53").Value = "Cultivar IV Región"
$ws.Range("I53").Value = "Especial"
$ws.Range("J53").Value = 50
$ws.Range("K53").Value = 20000
$ws.Range("L53").Value = 20000
$ws.Range("M53").Value = 20000
$ws.Range("N53").Value = "`$/bandeja 18 kilos"
$ws.Range("O53").Value = "Provincia de Limarí"
$ws.Range("P53").Value = 1111
$ws.Range("Q53").Value = 18
$ws.Range("R53").Value = "Hortaliza"

# New row 54: Primera quality entry also dated 2023-06-02 (serial 45079)
$ws.Range("A54").Value = 4
$ws.Range("B54").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C54").Value = "Los Lagos"
$ws.Range("D54").Value = 45079
$ws.Range("E54").Value = 10
$ws.Range("F54").Value = 100112043
$ws.Range("G54").Value = "Pepino dulce"
$ws.Range("H54").Value = "Cultivar IV Región"
$ws.Range("I54").Value = "Primera"
$ws.Range("J54").Value = 50
$ws.Range("K54").Value = 17000
$ws.Range("L54").Value = 17000
$ws.Range("M54").Value = 17000
$ws.Range("N54").Value = "`$/bandeja 18 kilos"
$ws.Range("O54").Value = "Provincia de Limarí"
$ws.Range("P54").Value = 944
$ws.Range("Q54").Value = 18
$ws.Range("R54").Value = "Hortaliza"
